$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '66.018.20'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -2.28%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.448.09'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  -0.03%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '584.62'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.10%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '173.93'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -3.16%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.81%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '3.446.65'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('E10').Value = '  -4.70%  '
$ws.Range('E11').Value = '  -1.26%  '
$ws.Range('E12').Value = '  -3.72%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '4.046.60'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -0.85%  '
$ws.Range('E14').Value = '  +1.12%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '29.04'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -9.81%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '66.048.59'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -2.26%  '
$ws.Range('E17').Value = '  -2.90%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '3.444.77'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -1.04%  '
$ws.Range('E19').Value = '  -3.09%  '
$ws.Range('E20').Value = '  -1.38%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '368.23'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -4.76%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '7.68'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -2.81%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '72.69'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('E24').Value = '  +0.07%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.537'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  +0.24%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.76'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('E28').Value = '  +0.83%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.08%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '23.79'
$c.Style = "Normal"
$ws.Range('E30').Value = '  -2.63%  '
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('E32').Value = '  -4.96%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').Value = '  -5.88%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '7.03'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -3.55%  '
$ws.Range('E36').Value = '  -1.75%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '161.08'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +0.52%  '
$ws.Range('E38').Value = '  +4.48%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.881'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -0.75%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '2.65'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -1.20%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '1.77'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -4.37%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '2.762.74'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +1.46%  '
$ws.Range('E43').Value = '  -1.57%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '6.42'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -3.42%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0682'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -3.80%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '40.16'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -3.21%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '24.32'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -5.75%  '
$ws.Range('E48').Value = '  -2.28%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '325.66'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.73%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '6.25'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.50%  '
$ws.Range('E51').Value = '  -2.83%  '
